# Refresh the coin Price (D) / Volume(1h) (E) columns to the values captured
# in the Mon Feb 26 18:50:26 UTC 2024 GitHub Actions run.
#
# Price cells that look like a plain decimal number (e.g. "1.00", "10.47")
# are written with a leading apostrophe. That is Excel's normal "force text"
# entry convention -- without it, Excel.Range.Value would silently reinterpret
# the digits as a numeric value (dropping the text formatting the source sheet
# uses for every Price/Volume cell, e.g. "1.00" -> 1). Percent strings and
# multi-dot price strings (e.g. "53.459.20") are never parsed as numbers, so
# they are written as plain text without the apostrophe.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '53.459.20'
$ws.Range("E2").Value = '  +3.72%  '

$ws.Range("D3").Value = '3.138.17'
$ws.Range("E3").Value = '  +2.56%  '

$ws.Range("E4").Value = '  +0.05%  '

$ws.Range("D5").Value = '''396.03'
$ws.Range("E5").Value = '  +2.47%  '

$ws.Range("D6").Value = '''109.45'
$ws.Range("E6").Value = '  +6.15%  '

$ws.Range("D8").Value = '''1.00'
$ws.Range("E8").Value = '  -0.02%  '

$ws.Range("D9").Value = '''0.609'
$ws.Range("E9").Value = '  +4.22%  '

$ws.Range("E10").Value = '  +5.33%  '

$ws.Range("E11").Value = '  +1.30%  '

$ws.Range("D12").Value = '''0.0871'
$ws.Range("E12").Value = '  +1.21%  '

$ws.Range("D13").Value = '3.643.00'
$ws.Range("E13").Value = '  +2.68%  '

$ws.Range("D14").Value = '''19.02'
$ws.Range("E14").Value = '  +2.23%  '

$ws.Range("E15").Value = '  +2.62%  '

$ws.Range("E16").Value = '  +8.04%  '

$ws.Range("D17").Value = '3.136.79'
$ws.Range("E17").Value = '  +2.64%  '

$ws.Range("D18").Value = '''10.47'
$ws.Range("E18").Value = '  -2.07%  '

$ws.Range("D19").Value = '53.357.73'
$ws.Range("E19").Value = '  +3.41%  '

$ws.Range("E20").Value = '  +3.30%  '

$ws.Range("D21").Value = '''12.72'
$ws.Range("E21").Value = '  +2.15%  '

$ws.Range("E22").Value = '  +0.50%  '

$ws.Range("D23").Value = '''70.84'
$ws.Range("E23").Value = '  +0.94%  '

$ws.Range("D24").Value = '''270.33'
$ws.Range("E24").Value = '  +0.87%  '

$ws.Range("D25").Value = '''3.23'
$ws.Range("E25").Value = '  +2.54%  '

$ws.Range("E26").Value = '  -2.82%  '

$ws.Range("D27").Value = '''27.42'
$ws.Range("E27").Value = '  +2.10%  '

$ws.Range("E28").Value = '  -0.39%  '

$ws.Range("E29").Value = '  -0.53%  '

$ws.Range("E30").Value = '  -0.09%  '

$ws.Range("E31").Value = '  +2.70%  '

$ws.Range("D32").Value = '''10.94'
$ws.Range("E32").Value = '  +6.40%  '

$ws.Range("E33").Value = '  +10.80%  '

$ws.Range("D34").Value = '''36.96'
$ws.Range("E34").Value = '  +6.19%  '

$ws.Range("E35").Value = '  +0.61%  '

$ws.Range("E36").Value = '  +0.90%  '

$ws.Range("D37").Value = '''3.65'
$ws.Range("E37").Value = '  +9.73%  '

$ws.Range("D38").Value = '''0.999'
$ws.Range("E38").Value = '  -0.12%  '

$ws.Range("D39").Value = '''2.77'
$ws.Range("E39").Value = '  +8.13%  '

$ws.Range("D40").Value = '''4.10'
$ws.Range("E40").Value = '  +9.41%  '

$ws.Range("E41").Value = '  -1.48%  '

$ws.Range("E42").Value = '  +1.60%  '

$ws.Range("E43").Value = '  +1.07%  '

$ws.Range("D44").Value = '''129.96'
$ws.Range("E44").Value = '  +3.91%  '

$ws.Range("E45").Value = '  +1.04%  '

$ws.Range("D46").Value = '''22.08'
$ws.Range("E46").Value = '  +0.58%  '

$ws.Range("E47").Value = '  -1.08%  '

$ws.Range("D48").Value = '''2.41'
$ws.Range("E48").Value = '  -0.26%  '

$ws.Range("D49").Value = '2.073.59'
$ws.Range("E49").Value = '  +1.92%  '

$ws.Range("D50").Value = '''0.0336'
$ws.Range("E50").Value = '  +5.53%  '

$ws.Range("D51").Value = '''0.0501'
$ws.Range("E51").Value = '  +15.81%  '
